$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 29.831634924547476;  C = 16.77035820586234;   D = 0.56216691603659175; E = 29.561514380455566; F = 14.988922401733872; G = 0.5070417641270677;  H = 277.92553153124362; I = 262.50856030781063 }
    3  = @{ B = 29.972035321394905;  C = 16.561317179833072;  D = 0.55255897713463342; E = 29.587399983523124; F = 14.702132219871187; G = 0.49690517680021334; H = 281.96892857101147; I = 264.48704707272367 }
    4  = @{ B = 30.560282556303289;  C = 16.532404964067737;  D = 0.54097683598340307; E = 29.556833147628854; F = 14.239183684070868; G = 0.48175606679341365; H = 288.41795033111913; I = 266 }
    5  = @{ B = 31.28384981176665;   C = 16.36691220972488;   D = 0.52317449125359461; E = 29.444516203785664; F = 14.223352846221534; G = 0.48305608921476678; H = 294;                I = 266.37946531598072 }
    6  = @{ B = 32.090165802622941;  C = 16.259890173305632;  D = 0.50669386606835798; E = 29.37831650694978;  F = 14.355390579492411; G = 0.48863897889099528; H = 298.51030057880655; I = 266.49719907237807 }
    7  = @{ B = 32.649814004987867;  C = 16.302775637972886;  D = 0.49932215955295589; E = 29.338398491322451; F = 14.519159892981312; G = 0.49488590514835734; H = 298.56552624829965; I = 266.06753678897422 }
    8  = @{ B = 33.068802040914271;  C = 16.123453276632159;  D = 0.48757294735634715; E = 29.3528167279946;   F = 14.561380840966121; G = 0.49608121005567846; H = 302.5;              I = 265.92052038816712 }
    9  = @{ B = 33.182609855255265;  C = 15.944136715839694;  D = 0.48049676578753359; E = 29.332295922129418; F = 14.614208585201702; G = 0.49822927683530488; H = 302.91498746270071; I = 265.90000000000009 }
    10 = @{ B = 33.151100490662976;  C = 15.80495628982203;   D = 0.47675510181852043; E = 29.336099060968252; F = 14.672402156332046; G = 0.50014837098275655; H = 302.5;              I = 266.05468782553243 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
